$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2021 column (R) -------------------------------------------------
# Header (row 4): year label, formatted like the existing year headers.
$ws.Range("R4").Value2 = 2021

# Data rows: renewable energy share (row 5) and hydropower production (row 6).
$ws.Range("R5").Value2 = 31.8
$ws.Range("R6").Value2 = 12957.1

# Copy the direct formatting from the neighbouring 2020 column (Q) / the
# column's own existing number formats so the new column matches the look
# of the rest of the table.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats (General, no border)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)   # xlPasteFormats (bottom border)

$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial of formats only shouldn't touch them,
# but make sure they are exactly right after formatting is copied over).
$ws.Range("R4").Value2 = 2021
$ws.Range("R5").Value2 = 31.8
$ws.Range("R6").Value2 = 12957.1

# Match the selection recorded in the saved workbook.
$ws.Range("R4:R6").Select() | Out-Null
